# This workbook's data rows (2-25, excluding row 23) were reshuffled:
# each destination row receives the full record (columns D and H..P)
# that originally belonged to a different source row, while columns
# A, B, C, E, F, G, Q, R stay identical (they were already constant for
# every row). Row 23 is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values are copied source -> destination)
$mapping = @{
    2  = 5
    3  = 16
    4  = 19
    5  = 4
    6  = 20
    7  = 17
    8  = 11
    9  = 12
    10 = 24
    11 = 13
    12 = 14
    13 = 3
    14 = 18
    15 = 7
    16 = 8
    17 = 21
    18 = 2
    19 = 9
    20 = 22
    21 = 15
    22 = 25
    24 = 6
    25 = 10
}

# Columns that vary per row and must be permuted (D, H..P).
$cols = @(4,8,9,10,11,12,13,14,15,16)

# Snapshot original values for every source row/column before overwriting
# anything, so that later writes don't clobber values still needed as a
# source for another destination row.
$original = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 25; $r++) {
        $original["$r-$col"] = $ws.Cells.Item($r, $col).Value2
    }
}

foreach ($dst in $mapping.Keys) {
    $src = $mapping[$dst]
    foreach ($col in $cols) {
        $ws.Cells.Item($dst, $col).Value2 = $original["$src-$col"]
    }
}
